$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1501.3334
$ws.Range("I28").Value = 1401.6154
$ws.Range("J28").Value = 2149.5
$ws.Range("K28").Value = 1401.6154
$ws.Range("L28").Value = 2149.5
$ws.Range("M28").Value = -916.6153999999999
$ws.Range("N28").Value = -3119.5

$ws.Range("H115").Value = 12000
$ws.Range("I115").Value = 9000
$ws.Range("K115").Value = 27000
$ws.Range("M115").Value = -25433

$ws.Range("H116").Value = 9147.637000000001
$ws.Range("I116").Value = 8788.4
$ws.Range("K116").Value = 8788.4
$ws.Range("M116").Value = -5346.4

$ws.Range("H118").Value = 2007.3334
$ws.Range("I118").Value = 1900
$ws.Range("J118").Value = 2222
$ws.Range("K118").Value = 5700
$ws.Range("L118").Value = 6666
$ws.Range("M118").Value = -4043
$ws.Range("N118").Value = -9980

$ws.Range("H125").Value = 761.6667
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").Value = $null

$ws.Range("H138").Value = 3232.9683
$ws.Range("J138").Value = 3298.7036
$ws.Range("L138").Value = 9896.110799999999
$ws.Range("N138").Value = -20176.1108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2498
$ws.Range("I2").Value = 2497
$ws.Range("J2").Value = 2498.5
$ws.Range("K2").Value = 2497
$ws.Range("L2").Value = 2498.5
$ws.Range("M2").Value = -2384
$ws.Range("N2").Value = -2724.5

$ws.Range("H116").Value = 2498
$ws.Range("I116").Value = 2497
$ws.Range("J116").Value = 2498.5
$ws.Range("K116").Value = 2497
$ws.Range("L116").Value = 2498.5
$ws.Range("M116").Value = -203
$ws.Range("N116").Value = -7086.5

$ws.Range("H122").Value = 25000
$ws.Range("I122").Value = 25000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 75000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -72550

$ws.Range("H132").Value = 3811.4285
$ws.Range("I132").Value = 2994.8333
$ws.Range("K132").Value = 8984.499899999999
$ws.Range("M132").Value = -6454.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2498
$ws.Range("I3").Value = 2497
$ws.Range("J3").Value = 2498.5
$ws.Range("K3").Value = 2497
$ws.Range("L3").Value = 2498.5
$ws.Range("M3").Value = -2383
$ws.Range("N3").Value = -2726.5

$ws.Range("H20").Value = 2354.375
$ws.Range("I20").Value = 691
$ws.Range("K20").Value = 691
$ws.Range("M20").Value = -444

$ws.Range("H94").Value = 1465.625
$ws.Range("I94").Value = 843.1667
$ws.Range("K94").Value = 843.1667
$ws.Range("M94").Value = -392.1667

$ws.Range("H96").Value = 32500
$ws.Range("I96").Value = 30000
$ws.Range("J96").Value = 35000
$ws.Range("K96").Value = 30000
$ws.Range("L96").Value = 35000
$ws.Range("M96").Value = -27254
$ws.Range("N96").Value = -40492

$ws.Range("H99").Value = 1100
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null

$ws.Range("H105").Value = 5616.25
$ws.Range("I105").Value = 4988.3335
$ws.Range("K105").Value = 4988.3335
$ws.Range("M105").Value = -3241.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1440.8
$ws.Range("I31").Value = 1402.5
$ws.Range("J31").Value = 1594
$ws.Range("K31").Value = 1402.5
$ws.Range("L31").Value = 1594
$ws.Range("M31").Value = -1107.5
$ws.Range("N31").Value = -2184

$ws.Range("H34").Value = 1440.8
$ws.Range("I34").Value = 1402.5
$ws.Range("J34").Value = 1594
$ws.Range("K34").Value = 1402.5
$ws.Range("L34").Value = 1594
$ws.Range("M34").Value = -1200.5
$ws.Range("N34").Value = -1998

$ws.Range("H58").Value = 2951.5833
$ws.Range("I58").Value = 2899
$ws.Range("J58").Value = 2962.1
$ws.Range("K58").Value = 2899
$ws.Range("L58").Value = 2962.1
$ws.Range("M58").Value = -2696
$ws.Range("N58").Value = -3368.1

$ws.Range("H86").Value = 17811.334
$ws.Range("I86").Value = 8636.916999999999
$ws.Range("K86").Value = 8636.916999999999
$ws.Range("M86").Value = -7513.916999999999

$ws.Range("H89").Value = 17811.334
$ws.Range("I89").Value = 8636.916999999999
$ws.Range("K89").Value = 43184.585
$ws.Range("M89").Value = -37568.585

$ws.Range("H100").Value = 67755
$ws.Range("J100").Value = 67755
$ws.Range("L100").Value = 67755
$ws.Range("N100").Value = -69919

$ws.Range("H107").Value = 1715.8948
$ws.Range("I107").Value = 745
$ws.Range("J107").Value = 2164
$ws.Range("K107").Value = 745
$ws.Range("L107").Value = 2164
$ws.Range("M107").Value = 1175
$ws.Range("N107").Value = -6004

$ws.Range("H122").Value = 4060.3
$ws.Range("I122").Value = 3772
$ws.Range("J122").Value = 4252.5
$ws.Range("K122").Value = 11316
$ws.Range("L122").Value = 12757.5
$ws.Range("M122").Value = -8866
$ws.Range("N122").Value = -17657.5

$ws.Range("H132").Value = 3477.1875
$ws.Range("I132").Value = 2879.3333
$ws.Range("J132").Value = 3835.9
$ws.Range("K132").Value = 8637.999899999999
$ws.Range("L132").Value = 11507.7
$ws.Range("M132").Value = -6107.999899999999
$ws.Range("N132").Value = -16567.7

$ws.Range("H136").Value = 2951.5833
$ws.Range("I136").Value = 2899
$ws.Range("J136").Value = 2962.1
$ws.Range("K136").Value = 8697
$ws.Range("L136").Value = 8886.299999999999
$ws.Range("M136").Value = -6147
$ws.Range("N136").Value = -13986.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 6692.6665
$ws.Range("I29").Value = 80
$ws.Range("J29").Value = 9999
$ws.Range("K29").Value = 240
$ws.Range("L29").Value = 29997
$ws.Range("M29").Value = 37
$ws.Range("N29").Value = -30551

$ws.Range("H40").Value = 200
$ws.Range("I40").Value = 200
$ws.Range("K40").Value = 800
$ws.Range("M40").Value = -731

$ws.Range("H80").Value = 204198.6
$ws.Range("I80").Value = 4497
$ws.Range("J80").Value = 337333
$ws.Range("K80").Value = 13491
$ws.Range("L80").Value = 1011999
$ws.Range("M80").Value = -12555
$ws.Range("N80").Value = -1013871

$ws.Range("H83").Value = 204198.6
$ws.Range("I83").Value = 4497
$ws.Range("J83").Value = 337333
$ws.Range("K83").Value = 40473
$ws.Range("L83").Value = 3035997
$ws.Range("M83").Value = -35793
$ws.Range("N83").Value = -3045357

$ws.Range("H122").Value = 101647.8
$ws.Range("J122").Value = 201799.4
$ws.Range("L122").Value = 1816194.6
$ws.Range("N122").Value = -1821094.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null

$ws.Range("H70").Value = 4721.6665
$ws.Range("J70").Value = 5285
$ws.Range("L70").Value = 5285
$ws.Range("N70").Value = -5825

$ws.Range("H73").Value = 4721.6665
$ws.Range("J73").Value = 5285
$ws.Range("L73").Value = 5285
$ws.Range("N73").Value = -7157

$ws.Range("H80").Value = 4497.1665
$ws.Range("I80").Value = 3395.4
$ws.Range("K80").Value = 3395.4
$ws.Range("M80").Value = -2397.4

$ws.Range("H83").Value = 4497.1665
$ws.Range("I83").Value = 3395.4
$ws.Range("K83").Value = 16977
$ws.Range("M83").Value = -11985

$ws.Range("H97").Value = 494.25
$ws.Range("I97").Value = 531.6429000000001
$ws.Range("K97").Value = 531.6429000000001
$ws.Range("M97").Value = -35.64290000000005

$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -6244

$ws.Range("H113").Value = 1367.4286
$ws.Range("I113").Value = 1367.4286
$ws.Range("K113").Value = 1367.4286
$ws.Range("M113").Value = 802.5714

$ws.Range("H132").Value = 3666
$ws.Range("I132").Value = 2999.4
$ws.Range("K132").Value = 8998.200000000001
$ws.Range("M132").Value = -6468.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2500
$ws.Range("I31").Value = 2500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -2252

$ws.Range("H46").Value = 1950
$ws.Range("I46").Value = 1933.3334
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 1933.3334
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = -1745.3334
$ws.Range("N46").Value = -2376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2652.2
$ws.Range("I136").Value = 2275
$ws.Range("J136").Value = 3322.7778
$ws.Range("K136").Value = 6825
$ws.Range("L136").Value = 9968.3334
$ws.Range("M136").Value = -4275
$ws.Range("N136").Value = -15068.3334
